$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.973193999999999
$ws.Range("H2").Value = 26.919582
$ws.Range("I2").Value = 0.3566934323398233
$ws.Range("J2").Value = 0.3566934323398234
$ws.Range("M2").Value = 5.904062666666666
$ws.Range("N2").Value = 17.712188
$ws.Range("O2").Value = 0.7745842164222705
$ws.Range("P2").Value = 0.7745842164222705
$ws.Range("Q2").Value = 52.97829969615732
$ws.Range("R2").Value = 476.8046972654159
$ws.Range("S2").Value = 0.2762891027919122
$ws.Range("T2").Value = 0.2762891027919123
# Row 3
$ws.Range("G3").Value = 8.973193999999999
$ws.Range("H3").Value = 26.919582
$ws.Range("I3").Value = 0.3566934323398233
$ws.Range("J3").Value = 0.3566934323398234
$ws.Range("M3").Value = 1.718172
$ws.Range("N3").Value = 5.154515999999999
$ws.Range("O3").Value = 0.2254157835777295
$ws.Range("P3").Value = 0.2254157835777294
$ws.Range("Q3").Value = 15.417490681368
$ws.Range("R3").Value = 138.757416132312
$ws.Range("S3").Value = 0.0804043295479111
$ws.Range("T3").Value = 0.0804043295479111
# Row 4
$ws.Range("I4").Value = 0.2126738244966861
$ws.Range("J4").Value = 0.2126738244966861
$ws.Range("M4").Value = 5.904062666666666
$ws.Range("N4").Value = 17.712188
$ws.Range("O4").Value = 0.7745842164222705
$ws.Range("P4").Value = 0.7745842164222705
$ws.Range("Q4").Value = 31.58762284408754
$ws.Range("R4").Value = 284.2886055967879
$ws.Range("S4").Value = 0.164733787701293
$ws.Range("T4").Value = 0.1647337877012931
# Row 5
$ws.Range("I5").Value = 0.2126738244966861
$ws.Range("J5").Value = 0.2126738244966861
$ws.Range("M5").Value = 1.718172
$ws.Range("N5").Value = 5.154515999999999
$ws.Range("O5").Value = 0.2254157835777295
$ws.Range("P5").Value = 0.2254157835777294
$ws.Range("Q5").Value = 9.192478498523998
$ws.Range("R5").Value = 82.73230648671598
$ws.Range("S5").Value = 0.047940036795393
$ws.Range("T5").Value = 0.04794003679539301
# Row 6
$ws.Range("G6").Value = 5.356989333333334
$ws.Range("H6").Value = 16.070968
$ws.Range("I6").Value = 0.2129456815838918
$ws.Range("J6").Value = 0.2129456815838919
$ws.Range("M6").Value = 5.904062666666666
$ws.Range("N6").Value = 17.712188
$ws.Range("O6").Value = 0.7745842164222705
$ws.Range("P6").Value = 0.7745842164222705
$ws.Range("Q6").Value = 31.62800072866489
$ws.Range("R6").Value = 284.652006557984
$ws.Range("S6").Value = 0.1649443639101652
$ws.Range("T6").Value = 0.1649443639101652
# Row 7
$ws.Range("G7").Value = 5.356989333333334
$ws.Range("H7").Value = 16.070968
$ws.Range("I7").Value = 0.2129456815838918
$ws.Range("J7").Value = 0.2129456815838919
$ws.Range("M7").Value = 1.718172
$ws.Range("N7").Value = 5.154515999999999
$ws.Range("O7").Value = 0.2254157835777295
$ws.Range("P7").Value = 0.2254157835777294
$ws.Range("Q7").Value = 9.204229076832
$ws.Range("R7").Value = 82.83806169148799
$ws.Range("S7").Value = 0.04800131767372665
$ws.Range("T7").Value = 0.04800131767372665
# Row 8
$ws.Range("G8").Value = 1.055528333333333
$ws.Range("H8").Value = 3.166585
$ws.Range("I8").Value = 0.04195830650140851
$ws.Range("J8").Value = 0.04195830650140851
$ws.Range("M8").Value = 5.904062666666666
$ws.Range("N8").Value = 17.712188
$ws.Range("O8").Value = 0.7745842164222705
$ws.Range("P8").Value = 0.7745842164222705
$ws.Range("Q8").Value = 6.231905426442221
$ws.Range("R8").Value = 56.08714883798
$ws.Range("S8").Value = 0.03250024196379897
$ws.Range("T8").Value = 0.03250024196379897
# Row 9
$ws.Range("G9").Value = 1.055528333333333
$ws.Range("H9").Value = 3.166585
$ws.Range("I9").Value = 0.04195830650140851
$ws.Range("J9").Value = 0.04195830650140851
$ws.Range("M9").Value = 1.718172
$ws.Range("N9").Value = 5.154515999999999
$ws.Range("O9").Value = 0.2254157835777295
$ws.Range("P9").Value = 0.2254157835777294
$ws.Range("Q9").Value = 1.81357922754
$ws.Range("R9").Value = 16.32221304786
$ws.Range("S9").Value = 0.009458064537609541
$ws.Range("T9").Value = 0.009458064537609539
# Row 10
$ws.Range("G10").Value = 4.420738
$ws.Range("H10").Value = 13.262214
$ws.Range("I10").Value = 0.1757287550781902
$ws.Range("J10").Value = 0.1757287550781902
$ws.Range("M10").Value = 5.904062666666666
$ws.Range("N10").Value = 17.712188
$ws.Range("O10").Value = 0.7745842164222705
$ws.Range("P10").Value = 0.7745842164222705
$ws.Range("Q10").Value = 26.10031418491466
$ws.Range("R10").Value = 234.902827664232
$ws.Range("S10").Value = 0.1361167200551011
$ws.Range("T10").Value = 0.1361167200551011
# Row 11
$ws.Range("G11").Value = 4.420738
$ws.Range("H11").Value = 13.262214
$ws.Range("I11").Value = 0.1757287550781902
$ws.Range("J11").Value = 0.1757287550781902
$ws.Range("M11").Value = 1.718172
$ws.Range("N11").Value = 5.154515999999999
$ws.Range("O11").Value = 0.2254157835777295
$ws.Range("P11").Value = 0.2254157835777294
$ws.Range("Q11").Value = 7.595588250935999
$ws.Range("R11").Value = 68.36029425842399
$ws.Range("S11").Value = 0.03961203502308915
$ws.Range("T11").Value = 0.03961203502308915
